$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(12,  "DK 5647 KM", "Friday, July 2, 2021", "Tidak Member", "Full Service", "Wax",        "Rp. 150,000", 200000, "Rp. 50,000"),
    @(15,  "DK 5678 HG", "Friday, July 2, 2021", "Tidak Member", "Half Service", "Wax, Engine","Rp. 230,000", 300000, "Rp. 70,000"),
    @(123, 123,          "Friday, July 2, 2021", "Tidak Member", "Full Service", "",           "Rp. 80,000",  200000, "Rp. 120,000")
)

$startRow = 12
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $col = $c + 1
        $val = $rowData[$c]
        if ($val -ne "") {
            $ws.Cells.Item($row, $col).Value = $val
        }
    }
}
